$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$green = 5296274  # RGB(146, 208, 80) -> FF92D050

# --- Row 2: add a "Java/Python" language tag for the Rotate-List question ---
$ws.Range("C2").Value = "Java/Python"

# --- New row 35: GFG doubly-linked-list pair-sum problem ---
$ws.Range("A35").Value = "GFG"
$ws.Range("B35").Value = "Find pairs with given sum in doubly linked list"
$ws.Range("C35").Value = "Java"

# --- Highlight (green fill) the rows that are now "done"/marked ---
$ws.Range("A2:C2").Interior.Color = $green
$ws.Range("A4:C4").Interior.Color = $green
$ws.Range("A5:B5").Interior.Color = $green
$ws.Range("A11:B11").Interior.Color = $green
$ws.Range("A12:C12").Interior.Color = $green
$ws.Range("A16:C16").Interior.Color = $green
$ws.Range("A17:C17").Interior.Color = $green
$ws.Range("A22:C22").Interior.Color = $green
$ws.Range("A32:C32").Interior.Color = $green
$ws.Range("A33:C33").Interior.Color = $green
$ws.Range("A34:C34").Interior.Color = $green
$ws.Range("A35:C35").Interior.Color = $green

# New row 35's B/C cells pick up stray wrap formatting inherited from the
# fill operation above (they sit past the sheet's previous used range) -
# force them back to the plain (non-wrapping) look used elsewhere in column C.
$ws.Range("B35:C35").WrapText = $false

# --- Update the view: scroll down a bit and leave the selection near the bottom ---
$ws.Range("B37").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 19
$win.ScrollColumn = 1
